$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "49.949.50"
Set-TextCell "E2" "  +3.98%  "
Set-TextCell "D3" "2.656.30"
Set-TextCell "E3" "  +6.21%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  +0.02%  "
Set-TextCell "D5" "328.07"
Set-TextCell "E5" "  +2.37%  "
Set-TextCell "D6" "111.29"
Set-TextCell "E6" "  +3.56%  "
Set-TextCell "E7" "  +1.27%  "
Set-TextCell "E8" "  +0.04%  "
Set-TextCell "D9" "0.560"
Set-TextCell "E9" "  +3.68%  "
Set-TextCell "D10" "40.83"
Set-TextCell "E10" "  +3.02%  "
Set-TextCell "D11" "20.56"
Set-TextCell "E11" "  +2.39%  "
Set-TextCell "E12" "  +1.50%  "
Set-TextCell "E13" "  +0.78%  "
Set-TextCell "D14" "7.34"
Set-TextCell "E14" "  +3.42%  "
Set-TextCell "D15" "3.071.24"
Set-TextCell "E15" "  +6.18%  "
Set-TextCell "D16" "2.638.34"
Set-TextCell "E16" "  +5.64%  "
Set-TextCell "D17" "0.882"
Set-TextCell "E17" "  +5.97%  "
Set-TextCell "D18" "49.930.91"
Set-TextCell "E18" "  +4.18%  "
Set-TextCell "E19" "  +2.82%  "
Set-TextCell "D20" "2.97"
Set-TextCell "E20" "  +7.72%  "
Set-TextCell "D21" "6.83"
Set-TextCell "E21" "  +1.77%  "
Set-TextCell "E22" "  +2.61%  "
Set-TextCell "D23" "73.18"
Set-TextCell "E23" "  +2.43%  "
Set-TextCell "D24" "280.99"
Set-TextCell "E24" "  +1.46%  "
Set-TextCell "E25" "  +2.78%  "
Set-TextCell "E26" "  +4.61%  "
Set-TextCell "D27" "0.999"
Set-TextCell "E27" "  -0.05%  "
Set-TextCell "B28" "InjectiveProtocol"
Set-TextCell "C28" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D28" "37.02"
Set-TextCell "E28" "  +5.07%  "
Set-TextCell "B29" "Toncoin"
Set-TextCell "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D29" "2.24"
Set-TextCell "E29" "  +7.44%  "
Set-TextCell "B30" "Cosmos"
Set-TextCell "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D30" "10.01"
Set-TextCell "E30" "  +3.08%  "
Set-TextCell "E31" "  +2.05%  "
Set-TextCell "D32" "49.78"
Set-TextCell "E32" "  +0.48%  "
Set-TextCell "B33" "Celestia"
Set-TextCell "C33" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D33" "19.63"
Set-TextCell "E33" "  +0.68%  "
Set-TextCell "B34" "Filecoin"
Set-TextCell "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D34" "5.45"
Set-TextCell "E34" "  +2.91%  "
Set-TextCell "E35" "  -0.10%  "
Set-TextCell "D36" "0.0801"
Set-TextCell "E36" "  +2.29%  "
Set-TextCell "E37" "  +7.18%  "
Set-TextCell "D38" "4.79"
Set-TextCell "E38" "  +3.68%  "
Set-TextCell "D39" "3.13"
Set-TextCell "E39" "  +8.99%  "
Set-TextCell "D40" "127.42"
Set-TextCell "E40" "  +5.34%  "
Set-TextCell "D42" "22.45"
Set-TextCell "E42" "  +6.10%  "
Set-TextCell "E43" "  +1.02%  "
Set-TextCell "E44" "  +3.86%  "
Set-TextCell "D45" "3.39"
Set-TextCell "E45" "  +7.82%  "
Set-TextCell "D46" "2.070.85"
Set-TextCell "E46" "  +2.62%  "
Set-TextCell "E47" "  +14.06%  "
Set-TextCell "E48" "  +8.42%  "
Set-TextCell "E49" "  +1.12%  "
Set-TextCell "E50" "  +5.04%  "
Set-TextCell "D51" "81.78"
Set-TextCell "E51" "  +1.69%  "
